# Add test for setFreezePanesTopCell:
# two more rows of data below the existing frozen header row, and move
# the freeze pane's anchor down so the new rows scroll in under it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (extends the used range from A1:A2 to A1:A4).
$ws.Range("A3").Value = "Should be at top after hello row when opening"
$ws.Range("A4").Value = "Another row"

# Re-apply the freeze so its pane reflects the new layout: still a
# 1-row freeze (ySplit=1, the "Hello" header stays pinned), but the
# scrollable pane's anchor moves on to A3 now that the sheet has grown.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
